$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting existing rows 240:250 down to 241:251
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new record
$ws.Range("A240").Value = 5
$ws.Range("B240").Value = "Macroferia Regional de Talca"
$ws.Range("C240").Value = "Maule"
$ws.Range("D240").Value = 44706
$ws.Range("E240").Value = 7
$ws.Range("F240").Value = 100112045
$ws.Range("G240").Value = "Zapallo"
$ws.Range("H240").Value = "Camote"
$ws.Range("I240").Value = "1a (guarda)"
$ws.Range("J240").Value = 900
$ws.Range("K240").Value = 400
$ws.Range("L240").Value = 400
$ws.Range("M240").Value = 400
$ws.Range("N240").Value = "$/kilo (volumen en unidades)"
$ws.Range("O240").Value = "Región del Maule"
$ws.Range("P240").Value = 400
$ws.Range("Q240").Value = 1
$ws.Range("R240").Value = "Hortaliza"
